$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data so the table can be rebuilt in its new row layout
$ws.UsedRange.ClearContents()

# Row 1 - new header row (added above the table)
$ws.Range('A1').Value = 'Next'
$ws.Range('B1').Value = 'Marlbro'
$ws.Range('C1').Value = 'LM'
$ws.Range('D1').Value = 'Galouis'
$ws.Range('E1').Value = 'Canadian'

# Row 2 - CO
$ws.Range('A2').Value = 'CO'
$ws.Range('B2').Value = 17.6394
$ws.Range('C2').Value = 6.5842
$ws.Range('D2').Value = 17.81
$ws.Range('E2').Value = 21.9388
$ws.Range('F2').Value = 21.9388

# Row 3 - CO2
$ws.Range('A3').Value = 'CO2'
$ws.Range('B3').Value = 177.0867
$ws.Range('C3').Value = 69.2488
$ws.Range('D3').Value = 174.368
$ws.Range('E3').Value = 239.1358
$ws.Range('F3').Value = 239.1358

# Row 4 - HCN
$ws.Range('A4').Value = 'HCN'
$ws.Range('D4').Value = 0.3698
$ws.Range('E4').Value = 0.5113
$ws.Range('F4').Value = 0.5113

# Row 5 - Methanol
$ws.Range('A5').Value = 'Methanol'
$ws.Range('D5').Value = 4.4847
$ws.Range('E5').Value = 4.2875
$ws.Range('F5').Value = 4.2875

# Row 6 - Ethylene
$ws.Range('A6').Value = 'Ethylene'
$ws.Range('D6').Value = 3.0965
$ws.Range('E6').Value = 3.1672
$ws.Range('F6').Value = 3.1672

# Row 7 - Isoprene
$ws.Range('A7').Value = 'Isoprene'
$ws.Range('D7').Value = 0.2788
$ws.Range('E7').Value = 0.3527
$ws.Range('F7').Value = 0.3527

# Row 8 - Methane
$ws.Range('A8').Value = 'Methane'
$ws.Range('D8').Value = 63.69
$ws.Range('E8').Value = 91.1204
$ws.Range('F8').Value = 91.1204

# Row 10 - Water
$ws.Range('A10').Value = 'Water'
$ws.Range('D10').Value = 43.1653
$ws.Range('E10').Value = 51.3689
$ws.Range('F10').Value = 51.3689

# Row 11 - Total
$ws.Range('A11').Value = 'Total'
$ws.Range('B11').Value = 194.7722
$ws.Range('C11').Value = 75.9313
$ws.Range('D11').Value = 523.1733
$ws.Range('E11').Value = 565.789
$ws.Range('F11').Value = 565.789

# Row 14 - Ethylene (% composition row; label unchanged from before)
$ws.Range('A14').Value = 'Ethylene'
$ws.Range('B14').Formula = '=C6/C$11 * 100'
$ws.Range('C14').Formula = '=D6/D$11 * 100'
$ws.Range('D14').Formula = '=E6/E$11 * 100'
$ws.Range('E14').Formula = '=F6/F$11 * 100'
$ws.Range('F14').Formula = '=B6/B$11 * 100'

# Row 15 - HCN (% composition row; label unchanged from before)
$ws.Range('A15').Value = 'HCN'
$ws.Range('B15').Formula = '=C4/C$11 * 100'
$ws.Range('C15').Formula = '=D4/D$11 * 100'
$ws.Range('D15').Formula = '=E4/E$11 * 100'
$ws.Range('E15').Formula = '=F4/F$11 * 100'
$ws.Range('F15').Formula = '=B4/B$11 * 100'

# Row 16 - Isoprene (% composition row; label unchanged from before)
$ws.Range('A16').Value = 'Isoprene'
$ws.Range('B16').Formula = '=C7/C$11 * 100'
$ws.Range('C16').Formula = '=D7/D$11 * 100'
$ws.Range('D16').Formula = '=E7/E$11 * 100'
$ws.Range('E16').Formula = '=F7/F$11 * 100'
$ws.Range('F16').Formula = '=B7/B$11 * 100'

# Row 17 - Methane (% composition row; label unchanged from before)
$ws.Range('A17').Value = 'Methane'
$ws.Range('B17').Formula = '=C8/C$11 * 100'
$ws.Range('C17').Formula = '=D8/D$11 * 100'
$ws.Range('D17').Formula = '=E8/E$11 * 100'
$ws.Range('E17').Formula = '=F8/F$11 * 100'
$ws.Range('F17').Formula = '=B8/B$11 * 100'

# Row 18 - Methanol (% composition row; label unchanged from before)
$ws.Range('A18').Value = 'Methanol'
$ws.Range('B18').Formula = '=C5/C$11 * 100'
$ws.Range('C18').Formula = '=D5/D$11 * 100'
$ws.Range('D18').Formula = '=E5/E$11 * 100'
$ws.Range('E18').Formula = '=F5/F$11 * 100'
$ws.Range('F18').Formula = '=B5/B$11 * 100'

# Row 19 - Water (% composition row; label unchanged from before)
$ws.Range('A19').Value = 'Water'
$ws.Range('B19').Formula = '=C10/C$11 * 100'
$ws.Range('C19').Formula = '=D10/D$11 * 100'
$ws.Range('D19').Formula = '=E10/E$11 * 100'
$ws.Range('E19').Formula = '=F10/F$11 * 100'
$ws.Range('F19').Formula = '=B10/B$11 * 100'

# Row 20 - CO (% composition row; label unchanged from before)
$ws.Range('A20').Value = 'CO'
$ws.Range('B20').Formula = '=C2/C$11 * 100'
$ws.Range('C20').Formula = '=D2/D$11 * 100'
$ws.Range('D20').Formula = '=E2/E$11 * 100'
$ws.Range('E20').Formula = '=F2/F$11 * 100'
$ws.Range('F20').Formula = '=B2/B$11 * 100'

# Row 21 - CO2 (% composition row; label unchanged from before)
$ws.Range('A21').Value = 'CO2'
$ws.Range('B21').Formula = '=C3/C$11 * 100'
$ws.Range('C21').Formula = '=D3/D$11 * 100'
$ws.Range('D21').Formula = '=E3/E$11 * 100'
$ws.Range('E21').Formula = '=F3/F$11 * 100'
$ws.Range('F21').Formula = '=B3/B$11 * 100'

# Row 22 - Sum (% composition row; label unchanged from before)
$ws.Range('A22').Value = 'Sum'
$ws.Range('B22').Formula = '=SUM(B14:B21)'
$ws.Range('C22').Formula = '=SUM(C14:C21)'
$ws.Range('D22').Formula = '=SUM(D14:D21)'
$ws.Range('E22').Formula = '=SUM(E14:E21)'
$ws.Range('F22').Formula = '=SUM(F14:F21)'

# Row 23 - Total (% composition row; label unchanged from before)
$ws.Range('A23').Value = 'Total'
$ws.Range('B23').Formula = '=C11/C$11 * 100'
$ws.Range('C23').Formula = '=D11/D$11 * 100'
$ws.Range('D23').Formula = '=E11/E$11 * 100'
$ws.Range('E23').Formula = '=F11/F$11 * 100'
$ws.Range('F23').Formula = '=B11/B$11%'

# Match the authored cursor/selection state
$ws.Range('E8').Select()
